$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new week's data (two rows: "Primera" and "Segunda" quality) is inserted at the
# top of the data block (row 487), pushing the existing rows 487:517 down to 489:519.
$ws.Rows("487:488").Insert()

# Row 487 - Primera
$ws.Cells.Item(487, 1).Value = 11
$ws.Cells.Item(487, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(487, 3).Value = "Bíobío"
$ws.Cells.Item(487, 4).Value = 45265
$ws.Cells.Item(487, 5).Value = 8
$ws.Cells.Item(487, 6).Value = 100112008
$ws.Cells.Item(487, 7).Value = "Coliflor"
$ws.Cells.Item(487, 8).Value = "Sin especificar"
$ws.Cells.Item(487, 9).Value = "Primera"
$ws.Cells.Item(487, 10).Value = 2000
$ws.Cells.Item(487, 11).Value = 900
$ws.Cells.Item(487, 12).Value = 1000
$ws.Cells.Item(487, 13).Value = 950
$ws.Cells.Item(487, 14).Value = "$/unidad"
$ws.Cells.Item(487, 15).Value = "Región Metropolitana"
$ws.Cells.Item(487, 16).Value = 950
$ws.Cells.Item(487, 17).Value = 1
$ws.Cells.Item(487, 18).Value = "Hortaliza"

# Row 488 - Segunda
$ws.Cells.Item(488, 1).Value = 11
$ws.Cells.Item(488, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(488, 3).Value = "Bíobío"
$ws.Cells.Item(488, 4).Value = 45265
$ws.Cells.Item(488, 5).Value = 8
$ws.Cells.Item(488, 6).Value = 100112008
$ws.Cells.Item(488, 7).Value = "Coliflor"
$ws.Cells.Item(488, 8).Value = "Sin especificar"
$ws.Cells.Item(488, 9).Value = "Segunda"
$ws.Cells.Item(488, 10).Value = 1000
$ws.Cells.Item(488, 11).Value = 800
$ws.Cells.Item(488, 12).Value = 800
$ws.Cells.Item(488, 13).Value = 800
$ws.Cells.Item(488, 14).Value = "$/unidad"
$ws.Cells.Item(488, 15).Value = "Región Metropolitana"
$ws.Cells.Item(488, 16).Value = 800
$ws.Cells.Item(488, 17).Value = 1
$ws.Cells.Item(488, 18).Value = "Hortaliza"
